$p = $ppt.ActivePresentation

# The deck currently ends with slide 21 ("Now code the StudentAssignments
# class yourself"). A new "Practice" slide needs to be inserted right
# before that final slide, i.e. at position 21 (pushing the old slide 21
# to become slide 22).
$lastIndex = $p.Slides.Count
$newSlide = $p.Slides.Add($lastIndex, 2)   # 2 = ppLayoutText ("Title and Content")

# Title placeholder
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Practice"

# Body / content placeholder - two paragraphs, second one indented one level
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Work in SmallClassProbs filling out ClassA.java, ClassB.java and ClassC.java `rshould work with the tests provided in ClassTests.java"

$para2 = $body.Paragraphs(2, 1)
$para2.IndentLevel = 2

# The old final slide (now pushed from position 21 to position 22) carries
# a notes page with a cached "slide number" field that still reads "21";
# refresh it to "22" now that the slide has moved.
$oldLastSlide = $p.Slides.Item($p.Slides.Count)
$notesPage = $oldLastSlide.NotesPage
for ($i = 1; $i -le $notesPage.Shapes.Count; $i++) {
    $shp = $notesPage.Shapes.Item($i)
    if ($shp.Name -eq "Slide Number Placeholder 3") {
        $shp.TextFrame.TextRange.Text = "22"
    }
}
